$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the user-entered answer values (D2:D7), leaving the formulas/formatting intact.
$ws.Range("D2:D7").ClearContents()

# Update the active selection to I8, matching the saved workbook view.
$ws.Range("I8").Select()
